# Update the "想去人数" (F column) counts that changed between scrapes.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 727
$wsExhibit.Range("F5").Value = 592
$wsExhibit.Range("F7").Value = 3458
$wsExhibit.Range("F8").Value = 481
$wsExhibit.Range("F9").Value = 8338
$wsExhibit.Range("F10").Value = 212
$wsExhibit.Range("F11").Value = 482
$wsExhibit.Range("F13").Value = 476

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 727
$wsAll.Range("F5").Value = 592
$wsAll.Range("F9").Value = 3458
$wsAll.Range("F10").Value = 481
$wsAll.Range("F12").Value = 8338
$wsAll.Range("F13").Value = 212
$wsAll.Range("F14").Value = 482
$wsAll.Range("F18").Value = 476
